$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 339, shifting existing rows 339:390 down to 340:391
$ws.Rows.Item(339).Insert()

# Populate the newly inserted row 339 with the new data entry.
# Columns A, B, C, E, F, G, H, I, R keep the same repeated values used throughout
# this data block; only D, J, K, L, M, N, O, P, Q differ for this record.
$ws.Cells.Item(339, 1).Value = 5
$ws.Cells.Item(339, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(339, 3).Value = "Maule"
$ws.Cells.Item(339, 4).Value = 44776
$ws.Cells.Item(339, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(339, 5).Value = 7
$ws.Cells.Item(339, 6).Value = 100112032
$ws.Cells.Item(339, 7).Value = "Zapallo italiano"
$ws.Cells.Item(339, 8).Value = "Sin especificar"
$ws.Cells.Item(339, 9).Value = "Primera"
$ws.Cells.Item(339, 10).Value = 200
$ws.Cells.Item(339, 11).Value = 20000
$ws.Cells.Item(339, 12).Value = 20000
$ws.Cells.Item(339, 13).Value = 20000
$ws.Cells.Item(339, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(339, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(339, 16).Value = 400
$ws.Cells.Item(339, 17).Value = 50
$ws.Cells.Item(339, 18).Value = "Hortaliza"
